$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user story rows to append (rows 26-28), mirroring style of rows 24-25
$stories = @(
    "AS A Buyer I should be able to Save the ADD-PRODUCT Details into Draft So THAT I can Continue the interrepted Comparison later. ",
    "AS A Buyer I Should be able to Clear the all Data which i was entered So THAT I can Start entering agin without any confusions. ",
    "AS A Buyer I should be able to Save the ADD-PRODUCT Details by iterations So THAT I can Continue the interrepted Comparison later. "
)

$startRow = 26
$slNo = 25

for ($i = 0; $i -lt $stories.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $slNo
    $ws.Cells.Item($row, 3).Value = $stories[$i]
    $ws.Cells.Item($row, 1).Style = $ws.Cells.Item(24, 1).Style
    $ws.Cells.Item($row, 3).Style = $ws.Cells.Item(24, 3).Style
    $ws.Rows.Item($row).RowHeight = $ws.Rows.Item(24).RowHeight
    $slNo++
}

$ws.Range("C29").Select()
